# Rows where Species (column F) is "Na": mark harvest as not occurring,
# clear the species cell, and zero out the unknown sex count.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,75,77,78,79)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = "No"
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 10).Value = 0
}
